# edit.ps1 -- applies the commit "Add files via upload" changes:
#   1. Remove the "-Vinit Prajapati (21BCP339)" and "-Rudra Shah (21BCP340)"
#      author subtitle paragraphs (and the blank paragraph that followed
#      them), so only "-Hrishikesh Kalola (21BCP346)" remains before the
#      "Introduction" heading.
#   2. Clean up the "Recommendation using GCN ... - https://www.mdpi.com/...”
#      citation hyperlink, which had its visible text split across three
#      runs ("http" / "s" / "://www.mdpi.com/2673-4591/58/1/97"); it should
#      be a single run reading "https://www.mdpi.com/2673-4591/58/1/97".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Delete the two extra author lines + the blank paragraph after them.
# ---------------------------------------------------------------------
# Locate the "-Vinit Prajapati (21BCP339)" paragraph and the
# "-Rudra Shah (21BCP340)" paragraph by scanning the Paragraphs
# collection (robust against any pre-existing paragraph-count
# assumptions), then remove the run of paragraphs from the first of
# those through the blank paragraph that immediately follows the
# second one.
$count = $d.Paragraphs.Count
$vinitIndex = 0
$rudraIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*Vinit Prajapati*") { $vinitIndex = $i }
    if ($ptext -like "*Rudra Shah*") { $rudraIndex = $i }
}

if ($vinitIndex -gt 0 -and $rudraIndex -gt $vinitIndex) {
    # The paragraph right after "Rudra Shah" is the stray blank
    # paragraph that should disappear together with the two name
    # lines (matches the source diff, which removes all three).
    $lastIndex = $rudraIndex
    $afterRudra = $d.Paragraphs.Item($rudraIndex + 1)
    if ($afterRudra.Range.Text.Trim() -eq "") {
        $lastIndex = $rudraIndex + 1
    }

    $startRange = $d.Paragraphs.Item($vinitIndex).Range
    $endRange = $d.Paragraphs.Item($lastIndex).Range
    $killRange = $d.Range($startRange.Start, $endRange.End)
    $killRange.Delete()
}

# ---------------------------------------------------------------------
# 2) Merge the split "http" / "s" / "://www.mdpi.com/..." hyperlink runs
#    into a single run, preserving the Hyperlink character style.
# ---------------------------------------------------------------------
$hCount = $d.Hyperlinks.Count
for ($i = 1; $i -le $hCount; $i++) {
    $link = $d.Hyperlinks.Item($i)
    $linkRange = $link.Range
    if ($linkRange.Text -eq "https://www.mdpi.com/2673-4591/58/1/97") {
        $firstRun = $linkRange.Duplicate
        $firstRun.Find.Execute("http", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

        $remainder = $d.Range($firstRun.End, $linkRange.End)
        $remainder.Delete()

        $onlyRun = $d.Range($firstRun.Start, $firstRun.End)
        $onlyRun.Find.Execute("http", $true, $false, $false, $false, $false, $true, 1, $false, "https://www.mdpi.com/2673-4591/58/1/97", 2) | Out-Null

        $fixedLink = $d.Hyperlinks.Item($i)
        $fixedLink.Range.Style = "Hyperlink"
        break
    }
}

Write-Output "done"
